$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Training Dashboard
$ws2 = $wb.Worksheets.Item(2)   # Exam Dashboard

# ---------------------------------------------------------------------
# 1) Header / title font recolor (white text).  The bold, size-14 title
#    font collapses onto the same bold/white font used by the header
#    band once the explicit 14pt size is dropped back to the default.
# ---------------------------------------------------------------------
foreach ($ws in @($ws1, $ws2)) {
    $ws.Range("A1").Font.Size = 11
    $ws.Range("A1").Font.Color = 16777215
}
$ws1.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215

# ---------------------------------------------------------------------
# 2) Training Dashboard (sheet1): refresh "PERIOD TO EXPIRE" (H) and
#    "LAST UPDATE" (I) for rows 3-28 -- last update moved from
#    08-Sep-2025 to 16-Sep-2025 (8 days later), so the day counts in H
#    all drop by 8.
# ---------------------------------------------------------------------
$hValues = @{
    3  = 399
    4  = 401
    5  = 495
    6  = 403
    7  = 410
    8  = 422
    9  = 520
    10 = 495
    11 = 520
    12 = 521
    13 = 521
    14 = 399
    15 = 420
    16 = 413
    17 = 523
    18 = 395
    19 = 522
    20 = 166
    21 = -102
    22 = -328
    23 = -49
    24 = 175
    25 = 215
    26 = 175
    27 = 217
    28 = 216
}

for ($row = 3; $row -le 28; $row++) {
    $ws1.Range("H$row").Value = $hValues[$row]

    # Write the date as literal text (not an auto-converted date serial)
    # by building it via a helper formula cell and pasting values only,
    # which preserves the original cell's existing number format/style.
    $ws1.Range("M1").Formula = '="16-Sep-2025"'
    $ws1.Range("M1").Copy() | Out-Null
    $ws1.Range("I$row").PasteSpecial(-4163) | Out-Null
}
$ws1.Range("M1").Clear() | Out-Null
$ws1.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Exam Dashboard (sheet2): comments column now reads "date is valid"
#    for the first 8 rows, and the COMMENTS column is narrower.
# ---------------------------------------------------------------------
$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"
$ws2.Range("E5").Value = "date is valid"
$ws2.Range("E6").Value = "date is valid"
$ws2.Range("E7").Value = "date is valid"
$ws2.Range("E8").Value = "date is valid"
$ws2.Range("E9").Value = "date is valid"
$ws2.Range("E10").Value = "date is valid"

$ws2.Columns.Item(5).ColumnWidth = 14.17
